{"js": "// 1) Insert a new \"Meta description\" paragraph right after the title (first) paragraph.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nconst metaPara = titlePara.insertParagraph(\"\", \"After\");\nmetaPara.style = \"Normal\";\nmetaPara.insertText(\n  \"Meta description: Discover the magical world of Alkemor's Tower, a unique and exciting slot game. Learn how to play and trigger its special functions for better winnings. Play for free.\",\n  \"End\"\n);\nawait context.sync();\n\n// Make just the \"Meta description\" portion bold (leaving the rest of the sentence regular).\nconst metaLabel = metaPara.search(\"Meta description\", { matchCase: true });\nmetaLabel.load(\"items\");\nawait context.sync();\nmetaLabel.items[0].font.bold = true;\nawait context.sync();\n\n// 2) Remove the duplicated bold title paragraph that was left near the end of the document,\n//    and 3) repurpose the italic paragraph that follows it into an image-generation prompt.\n// NOTE: the document's very first paragraph (the Heading 1 title) has the same text as the\n// duplicate paragraph near the end, so search from the end of the document to find the right one.\nconst bodyParagraphs = context.document.body.paragraphs;\nbodyParagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nlet duplicateTitle = null;\nfor (let i = bodyParagraphs.items.length - 1; i >= 0; i--) {\n  const p = bodyParagraphs.items[i];\n  if (\n    p.text === \"Play Alkemor's Tower for Free - Review & Gameplay Mechanics\" &&\n    p.style !== \"Heading 1\"\n  ) {\n    duplicateTitle = p;\n    break;\n  }\n}\nif (duplicateTitle) {\n  duplicateTitle.delete();\n  await context.sync();\n}\n\nconst oldDescriptionText =\n  \"Discover the magical world of Alkemor's Tower, a unique and exciting slot game. Learn how to play and trigger its special functions for better winnings. Play for free.\";\nconst newPromptText =\n  \"Please create a feature image for Alkemor's Tower that fits the following specifications: - Cartoon style - Features a happy Maya warrior with glasses\";\n\nconst descriptionParagraphs = context.document.body.paragraphs;\ndescriptionParagraphs.load(\"items/text\");\nawait context.sync();\n\nconst promptPara = descriptionParagraphs.items.find((p) => p.text === oldDescriptionText);\nif (promptPara) {\n  const found = promptPara.search(oldDescriptionText, { matchCase: true });\n  found.load(\"items\");\n  await context.sync();\n  found.items[0].insertText(newPromptText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1) Insert a new \"Meta description\" paragraph right after the title (first) paragraph. ---\n$nextPara = $d.Paragraphs.Item(2)\n$nextPara.Range.InsertParagraphBefore()\n\n# The freshly inserted paragraph is now #2; it inherited paragraph #3's formatting, so reset it.\n$metaPara = $d.Paragraphs.Item(2)\n$metaPara.Style = \"Normal\"\n$metaPara.Range.InsertAfter(\"Meta description: Discover the magical world of Alkemor's Tower, a unique and exciting slot game. Learn how to play and trigger its special functions for better winnings. Play for free.\")\n\n# Bold just the \"Meta description\" label, leaving the rest of the sentence regular.\n$labelRange = $metaPara.Range.Duplicate\n[void]$labelRange.Find.Execute(\"Meta description\")\n$labelRange.Bold = 1\n\n# --- 2) Remove the duplicated bold title paragraph left near the end of the document. ---\n# (The very first paragraph - the Heading 1 title - has identical text, so disambiguate by style\n#  and search from the end of the document.)\n$titleText = \"Play Alkemor's Tower for Free - Review & Gameplay Mechanics\"\n$count = $d.Paragraphs.Count\n$dupIndex = -1\nfor ($i = $count; $i -ge 1; $i--) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n  if ($t -eq $titleText -and $p.Style.NameLocal -ne \"Heading 1\") {\n    $dupIndex = $i\n    break\n  }\n}\nif ($dupIndex -ge 1) {\n  $d.Paragraphs.Item($dupIndex).Range.Delete()\n}\n\n# --- 3) Repurpose the italic paragraph that follows it into an image-generation prompt. ---\n$oldDescriptionText = \"Discover the magical world of Alkemor's Tower, a unique and exciting slot game. Learn how to play and trigger its special functions for better winnings. Play for free.\"\n$newPromptText = \"Please create a feature image for Alkemor's Tower that fits the following specifications: - Cartoon style - Features a happy Maya warrior with glasses\"\n\n$count2 = $d.Paragraphs.Count\n$descIndex = -1\nfor ($i = $count2; $i -ge 1; $i--) {\n  $p = $d.Paragraphs.Item($i)\n  $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n  if ($t -eq $oldDescriptionText) {\n    $descIndex = $i\n    break\n  }\n}\nif ($descIndex -ge 1) {\n  $descPara = $d.Paragraphs.Item($descIndex)\n  $findRange = $descPara.Range.Duplicate\n  [void]$findRange.Find.Execute($oldDescriptionText)\n  $findRange.Text = $newPromptText\n}\n"}
